# Auto-generated script to apply odds updates per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.2  # G2: 1.73 -> 2.2
$ws.Cells.Item(2, 8).Value = 2.9  # H2: 3.25 -> 2.9
$ws.Cells.Item(2, 9).Value = 3.7  # I2: 5.5 -> 3.7
$ws.Cells.Item(2, 10).Value = 3.1  # J2: 2.5 -> 3.1
$ws.Cells.Item(2, 11).Value = 1.91  # K2: 1.95 -> 1.91
$ws.Cells.Item(2, 12).Value = 4.5  # L2: 6 -> 4.5
$ws.Cells.Item(2, 13).Value = 1.11  # M2: 1.1 -> 1.11
$ws.Cells.Item(2, 14).Value = 6.5  # N2: 7 -> 6.5
$ws.Cells.Item(2, 15).Value = 1.53  # O2: 1.5 -> 1.53
$ws.Cells.Item(2, 16).Value = 2.38  # P2: 2.5 -> 2.38
$ws.Cells.Item(2, 17).Value = 2.05  # Q2: 1.98 -> 2.05
$ws.Cells.Item(2, 18).Value = 1.8  # R2: 1.88 -> 1.8
$ws.Cells.Item(2, 19).Value = 2.7  # S2: 2.6 -> 2.7
$ws.Cells.Item(2, 20).Value = 1.44  # T2: 1.48 -> 1.44
$ws.Cells.Item(2, 21).Value = 4.5  # U2: 4.1 -> 4.5
$ws.Cells.Item(2, 22).Value = 1.2  # V2: 1.22 -> 1.2
$ws.Cells.Item(2, 23).Value = 5.5  # W2: 5 -> 5.5
$ws.Cells.Item(2, 24).Value = 1.14  # X2: 1.17 -> 1.14
$ws.Cells.Item(2, 25).Value = 1.62  # Y2: 1.57 -> 1.62
$ws.Cells.Item(2, 26).Value = 2.2  # Z2: 2.25 -> 2.2
$ws.Cells.Item(2, 27).Value = 2.2  # AA2: 2.38 -> 2.2
$ws.Cells.Item(2, 28).Value = 1.62  # AB2: 1.53 -> 1.62
$ws.Cells.Item(2, 29).Value = 5.5  # AC2: 5 -> 5.5
$ws.Cells.Item(2, 30).Value = 9  # AD2: 6.5 -> 9
$ws.Cells.Item(2, 31).Value = 10  # AE2: 9.5 -> 10
$ws.Cells.Item(2, 32).Value = 21  # AF2: 13 -> 21
$ws.Cells.Item(2, 33).Value = 23  # AG2: 19 -> 23
$ws.Cells.Item(2, 35).Value = 6  # AI2: 6.5 -> 6
$ws.Cells.Item(2, 36).Value = 6  # AJ2: 6.5 -> 6
$ws.Cells.Item(2, 37).Value = 19  # AK2: 21 -> 19
$ws.Cells.Item(2, 40).Value = 8  # AN2: 11 -> 8
$ws.Cells.Item(2, 41).Value = 17  # AO2: 26 -> 17
$ws.Cells.Item(2, 42).Value = 15  # AP2: 19 -> 15
$ws.Cells.Item(2, 43).Value = 41  # AQ2: 67 -> 41
$ws.Cells.Item(2, 44).Value = 41  # AR2: 51 -> 41
$ws.Cells.Item(2, 45).Value = 51  # AS2: 67 -> 51
# Row 3
$ws.Cells.Item(3, 21).Value = 4  # U3: 4.1 -> 4
# Row 4
$ws.Cells.Item(4, 7).Value = 2.75  # G4: 2.9 -> 2.75
$ws.Cells.Item(4, 9).Value = 2.9  # I4: 2.75 -> 2.9
$ws.Cells.Item(4, 10).Value = 3.75  # J4: 4 -> 3.75
$ws.Cells.Item(4, 12).Value = 4  # L4: 3.75 -> 4
$ws.Cells.Item(4, 15).Value = 1.67  # O4: 1.62 -> 1.67
$ws.Cells.Item(4, 16).Value = 2.1  # P4: 2.2 -> 2.1
$ws.Cells.Item(4, 29).Value = 6  # AC4: 6.5 -> 6
$ws.Cells.Item(4, 30).Value = 11  # AD4: 12 -> 11
$ws.Cells.Item(4, 31).Value = 12  # AE4: 13 -> 12
$ws.Cells.Item(4, 32).Value = 29  # AF4: 34 -> 29
$ws.Cells.Item(4, 33).Value = 29  # AG4: 34 -> 29
$ws.Cells.Item(4, 40).Value = 6.5  # AN4: 6 -> 6.5
$ws.Cells.Item(4, 41).Value = 12  # AO4: 11 -> 12
$ws.Cells.Item(4, 42).Value = 13  # AP4: 12 -> 13
$ws.Cells.Item(4, 43).Value = 34  # AQ4: 29 -> 34
$ws.Cells.Item(4, 44).Value = 34  # AR4: 29 -> 34
# Row 5
$ws.Cells.Item(5, 12).Value = 3.65  # L5: 3.7 -> 3.65
$ws.Cells.Item(5, 14).Value = 4.15  # N5: 4.2 -> 4.15
$ws.Cells.Item(5, 15).Value = 1.75  # O5: 1.72 -> 1.75
$ws.Cells.Item(5, 16).Value = 1.98  # P5: 2 -> 1.98
$ws.Cells.Item(5, 25).Value = 1.7  # Y5: 1.72 -> 1.7
$ws.Cells.Item(5, 26).Value = 2.02  # Z5: 2 -> 2.02
$ws.Cells.Item(5, 27).Value = 2.3  # AA5: 2.32 -> 2.3
$ws.Cells.Item(5, 35).Value = 4.15  # AI5: 4.2 -> 4.15
$ws.Cells.Item(5, 40).Value = 5.7  # AN5: 5.6 -> 5.7
$ws.Cells.Item(5, 42).Value = 11.5  # AP5: 11.75 -> 11.5
# Row 6
$ws.Cells.Item(6, 7).Value = 2.95  # G6: 2.92 -> 2.95
$ws.Cells.Item(6, 8).Value = 2.62  # H6: 2.65 -> 2.62
$ws.Cells.Item(6, 10).Value = 3.75  # J6: 3.6 -> 3.75
$ws.Cells.Item(6, 11).Value = 1.8  # K6: 1.83 -> 1.8
$ws.Cells.Item(6, 12).Value = 3.45  # L6: 3.55 -> 3.45
$ws.Cells.Item(6, 13).Value = 1.17  # M6: 1.16 -> 1.17
$ws.Cells.Item(6, 14).Value = 4.45  # N6: 4.55 -> 4.45
$ws.Cells.Item(6, 19).Value = 3.05  # S6: 2.95 -> 3.05
$ws.Cells.Item(6, 20).Value = 1.33  # T6: 1.35 -> 1.33
$ws.Cells.Item(6, 23).Value = 5.6  # W6: 5.5 -> 5.6
$ws.Cells.Item(6, 26).Value = 2.15  # Z6: 2.18 -> 2.15
$ws.Cells.Item(6, 27).Value = 2.32  # AA6: 2.3 -> 2.32
$ws.Cells.Item(6, 28).Value = 1.53  # AB6: 1.55 -> 1.53
$ws.Cells.Item(6, 29).Value = 5.8  # AC6: 6.2 -> 5.8
$ws.Cells.Item(6, 30).Value = 12.5  # AD6: 13 -> 12.5
$ws.Cells.Item(6, 31).Value = 12.5  # AE6: 11.5 -> 12.5
$ws.Cells.Item(6, 32).Value = 40  # AF6: 37 -> 40
$ws.Cells.Item(6, 33).Value = 40  # AG6: 35 -> 40
$ws.Cells.Item(6, 34).Value = 70  # AH6: 60 -> 70
$ws.Cells.Item(6, 35).Value = 4.45  # AI6: 4.55 -> 4.45
$ws.Cells.Item(6, 37).Value = 21  # AK6: 20 -> 21
$ws.Cells.Item(6, 38).Value = 175  # AL6: 150 -> 175
$ws.Cells.Item(6, 40).Value = 5.8  # AN6: 5.7 -> 5.8
$ws.Cells.Item(6, 41).Value = 12  # AO6: 11.75 -> 12
$ws.Cells.Item(6, 42).Value = 11.25  # AP6: 11.75 -> 11.25
$ws.Cells.Item(6, 44).Value = 32  # AR6: 35 -> 32
$ws.Cells.Item(6, 45).Value = 60  # AS6: 65 -> 60
# Row 7
$ws.Cells.Item(7, 7).Value = 2.45  # G7: 2.4 -> 2.45
$ws.Cells.Item(7, 9).Value = 3.25  # I7: 3.3 -> 3.25
$ws.Cells.Item(7, 27).Value = 2.5  # AA7: 2.38 -> 2.5
$ws.Cells.Item(7, 28).Value = 1.5  # AB7: 1.53 -> 1.5
$ws.Cells.Item(7, 30).Value = 10  # AD7: 9.5 -> 10
$ws.Cells.Item(7, 32).Value = 26  # AF7: 23 -> 26
$ws.Cells.Item(7, 33).Value = 29  # AG7: 26 -> 29
$ws.Cells.Item(7, 41).Value = 13  # AO7: 15 -> 13
$ws.Cells.Item(7, 44).Value = 34  # AR7: 41 -> 34
# Row 9
$ws.Cells.Item(9, 7).Value = 2.35  # G9: 2.38 -> 2.35
$ws.Cells.Item(9, 8).Value = 3.2  # H9: 3.25 -> 3.2
$ws.Cells.Item(9, 9).Value = 3  # I9: 2.9 -> 3
$ws.Cells.Item(9, 11).Value = 1.91  # K9: 1.95 -> 1.91
$ws.Cells.Item(9, 12).Value = 4  # L9: 3.75 -> 4
$ws.Cells.Item(9, 13).Value = 1.1  # M9: 1.08 -> 1.1
$ws.Cells.Item(9, 14).Value = 7  # N9: 8 -> 7
$ws.Cells.Item(9, 15).Value = 1.5  # O9: 1.44 -> 1.5
$ws.Cells.Item(9, 16).Value = 2.5  # P9: 2.63 -> 2.5
$ws.Cells.Item(9, 17).Value = 1.9  # Q9: 1.85 -> 1.9
$ws.Cells.Item(9, 18).Value = 1.95  # R9: 2 -> 1.95
$ws.Cells.Item(9, 19).Value = 2.5  # S9: 2.4 -> 2.5
$ws.Cells.Item(9, 20).Value = 1.5  # T9: 1.53 -> 1.5
$ws.Cells.Item(9, 25).Value = 1.57  # Y9: 1.53 -> 1.57
$ws.Cells.Item(9, 26).Value = 2.25  # Z9: 2.38 -> 2.25
$ws.Cells.Item(9, 27).Value = 2.1  # AA9: 2.05 -> 2.1
$ws.Cells.Item(9, 28).Value = 1.67  # AB9: 1.7 -> 1.67
$ws.Cells.Item(9, 29).Value = 6  # AC9: 6.5 -> 6
$ws.Cells.Item(9, 35).Value = 7  # AI9: 7.5 -> 7
$ws.Cells.Item(9, 40).Value = 7  # AN9: 7.5 -> 7
$ws.Cells.Item(9, 42).Value = 12  # AP9: 11 -> 12
# Row 10
$ws.Cells.Item(10, 7).Value = 3.8  # G10: 3.75 -> 3.8
$ws.Cells.Item(10, 8).Value = 3.2  # H10: 3.25 -> 3.2
$ws.Cells.Item(10, 11).Value = 1.91  # K10: 1.95 -> 1.91
$ws.Cells.Item(10, 17).Value = 1.93  # Q10: 1.98 -> 1.93
$ws.Cells.Item(10, 18).Value = 1.93  # R10: 1.88 -> 1.93
$ws.Cells.Item(10, 19).Value = 2.5  # S10: 2.6 -> 2.5
$ws.Cells.Item(10, 20).Value = 1.5  # T10: 1.48 -> 1.5
$ws.Cells.Item(10, 23).Value = 5  # W10: 5.5 -> 5
$ws.Cells.Item(10, 24).Value = 1.17  # X10: 1.14 -> 1.17
$ws.Cells.Item(10, 27).Value = 2.2  # AA10: 2.1 -> 2.2
$ws.Cells.Item(10, 28).Value = 1.62  # AB10: 1.67 -> 1.62
$ws.Cells.Item(10, 29).Value = 8  # AC10: 8.5 -> 8
$ws.Cells.Item(10, 31).Value = 15  # AE10: 13 -> 15
$ws.Cells.Item(10, 34).Value = 51  # AH10: 41 -> 51
$ws.Cells.Item(10, 35).Value = 6.5  # AI10: 7 -> 6.5
$ws.Cells.Item(10, 37).Value = 21  # AK10: 19 -> 21
$ws.Cells.Item(10, 41).Value = 8  # AO10: 8.5 -> 8
# Row 11
$ws.Cells.Item(11, 13).Value = 1.1  # M11: 1.08 -> 1.1
$ws.Cells.Item(11, 14).Value = 7  # N11: 7.5 -> 7
$ws.Cells.Item(11, 15).Value = 1.44  # O11: 1.4 -> 1.44
$ws.Cells.Item(11, 16).Value = 2.63  # P11: 2.75 -> 2.63
$ws.Cells.Item(11, 35).Value = 7  # AI11: 7.5 -> 7
# Row 12
$ws.Cells.Item(12, 7).Value = 2.2  # G12: 2.25 -> 2.2
$ws.Cells.Item(12, 9).Value = 3.4  # I12: 3.25 -> 3.4
$ws.Cells.Item(12, 11).Value = 1.91  # K12: 1.95 -> 1.91
$ws.Cells.Item(12, 12).Value = 4.33  # L12: 4 -> 4.33
$ws.Cells.Item(12, 13).Value = 1.08  # M12: 1.1 -> 1.08
$ws.Cells.Item(12, 14).Value = 8  # N12: 7 -> 8
$ws.Cells.Item(12, 15).Value = 1.44  # O12: 1.5 -> 1.44
$ws.Cells.Item(12, 16).Value = 2.63  # P12: 2.5 -> 2.63
$ws.Cells.Item(12, 17).Value = 1.85  # Q12: 1.9 -> 1.85
$ws.Cells.Item(12, 18).Value = 2  # R12: 1.95 -> 2
$ws.Cells.Item(12, 19).Value = 2.4  # S12: 2.5 -> 2.4
$ws.Cells.Item(12, 20).Value = 1.53  # T12: 1.5 -> 1.53
$ws.Cells.Item(12, 27).Value = 2.1  # AA12: 2.05 -> 2.1
$ws.Cells.Item(12, 28).Value = 1.67  # AB12: 1.7 -> 1.67
$ws.Cells.Item(12, 35).Value = 6.5  # AI12: 7 -> 6.5
$ws.Cells.Item(12, 40).Value = 8  # AN12: 7.5 -> 8
$ws.Cells.Item(12, 42).Value = 13  # AP12: 12 -> 13
$ws.Cells.Item(12, 43).Value = 41  # AQ12: 34 -> 41
# Row 13
$ws.Cells.Item(13, 9).Value = 7.5  # I13: 8 -> 7.5
$ws.Cells.Item(13, 10).Value = 1.91  # J13: 1.95 -> 1.91
$ws.Cells.Item(13, 12).Value = 7  # L13: 7.5 -> 7
$ws.Cells.Item(13, 13).Value = 1.04  # M13: 1.05 -> 1.04
$ws.Cells.Item(13, 14).Value = 13  # N13: 11 -> 13
$ws.Cells.Item(13, 15).Value = 1.22  # O13: 1.25 -> 1.22
$ws.Cells.Item(13, 16).Value = 4  # P13: 3.75 -> 4
$ws.Cells.Item(13, 19).Value = 1.8  # S13: 1.88 -> 1.8
$ws.Cells.Item(13, 20).Value = 2  # T13: 1.98 -> 2
$ws.Cells.Item(13, 23).Value = 3  # W13: 3.25 -> 3
$ws.Cells.Item(13, 24).Value = 1.36  # X13: 1.33 -> 1.36
$ws.Cells.Item(13, 27).Value = 2  # AA13: 2.05 -> 2
$ws.Cells.Item(13, 28).Value = 1.75  # AB13: 1.7 -> 1.75
$ws.Cells.Item(13, 29).Value = 7  # AC13: 6.5 -> 7
$ws.Cells.Item(13, 39).Value = 351  # AM13: 401 -> 351
$ws.Cells.Item(13, 43).Value = 81  # AQ13: 101 -> 81
# Row 14
$ws.Cells.Item(14, 7).Value = 2.15  # G14: 2.1 -> 2.15
$ws.Cells.Item(14, 13).Value = 1.03  # M14: 1.06 -> 1.03
$ws.Cells.Item(14, 14).Value = 9.5  # N14: 8 -> 9.5
$ws.Cells.Item(14, 19).Value = 2.03  # S14: 2.05 -> 2.03
$ws.Cells.Item(14, 20).Value = 1.83  # T14: 1.8 -> 1.83
$ws.Cells.Item(14, 23).Value = 3.4  # W14: 3.5 -> 3.4
$ws.Cells.Item(14, 24).Value = 1.3  # X14: 1.29 -> 1.3

Write-Host "Applied 183 cell updates"
